# dif_semanal_Vfarnesiana_2.xlsx — "Add files via upload"
#
# The uploaded workbook is a re-save of the same weekly-difference dataset
# with 12 replicate rows removed (the underlying data for those replicas
# was dropped from the source table). No other cell values changed - once
# the rows are removed the remaining rows simply shift up, which is what
# made the raw XML diff look like a wall of value changes.
#
# Identify the rows to drop by their "Réplica" id (column A) rather than by
# a hard-coded worksheet row number, so the deletion is robust/self-describing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Réplica ids (column A) that must disappear from the table.
$replicaIdsToRemove = @(3, 29, 31, 37, 40, 41, 44, 51, 53, 54, 55, 66)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Walk bottom-to-top so deleting a row never invalidates the row numbers of
# rows we still need to inspect/delete above it.
for ($r = $lastRow; $r -ge 2; $r--) {
    $replicaId = $ws.Cells.Item($r, 1).Value2
    if ($replicaIdsToRemove -contains $replicaId) {
        $ws.Rows($r).Delete()
    }
}

# Return the selection to the top-left cell (the stale A1:D51 selection from
# before the deletions no longer matches the shrunk A1:D39 used range).
$ws.Range("A1").Select()

Write-Host "Removed rows for replicas:" ($replicaIdsToRemove -join ", ")
